# Trading update: 2026-02-17 15:32:26
# Appends the newest (still OPEN) MarketMaking trade as row 50 to both the
# "All Trades" and "MarketMaking" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 50

    # Text columns: force text format first so date/time-like strings
    # ("2026-02-17") are not auto-coerced into date serial numbers.
    $ws.Cells.Item($row, 2).NumberFormat = "@"   # Date
    $ws.Cells.Item($row, 3).NumberFormat = "@"   # Time
    $ws.Cells.Item($row, 4).NumberFormat = "@"   # Strategy
    $ws.Cells.Item($row, 5).NumberFormat = "@"   # Side
    $ws.Cells.Item($row, 8).NumberFormat = "@"   # Status
    $ws.Cells.Item($row, 15).NumberFormat = "@"  # Entry Reason
    $ws.Cells.Item($row, 16).NumberFormat = "@"  # Exit Reason

    $ws.Cells.Item($row, 1).Value = 49                                        # A: Trade #
    $ws.Cells.Item($row, 2).Value = "2026-02-17"                              # B: Date
    $ws.Cells.Item($row, 3).Value = "15:30:58"                                # C: Time
    $ws.Cells.Item($row, 4).Value = "MarketMaking"                            # D: Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"                                    # E: Side
    $ws.Cells.Item($row, 6).Value = 0.49                                      # F: Entry Price
    $ws.Cells.Item($row, 7).Value = ""                                        # G: Exit Price (empty, still open)
    $ws.Cells.Item($row, 8).Value = "OPEN"                                    # H: Status
    $ws.Cells.Item($row, 9).Value = 0                                         # I: P&L %
    $ws.Cells.Item($row, 10).Value = 0                                        # J: P&L $
    $ws.Cells.Item($row, 11).Value = 100.5515569553527                        # K: Capital After
    $ws.Cells.Item($row, 12).Value = 0                                        # L: Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                                        # M: Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                                      # N: Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"       # O: Entry Reason
    $ws.Cells.Item($row, 16).Value = ""                                       # P: Exit Reason (empty, still open)
    $ws.Cells.Item($row, 17).Value = 0                                        # Q: Duration (min)
}
